# updatedDBs and added issue #22
#
# Appends 7 new Korean-vocabulary entries (FOREIGN / ENGLISH / DATE) to the
# end of the "words" sheet, continuing on from the existing last row (76).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newWords = @(
    @{ Foreign = "분야";       English = "1.)  n. sphere, area, field, domain; 2.) n. branch, realm, province" },
    @{ Foreign = "진출하다";   English = "1.) v. to branch out, expand" },
    @{ Foreign = "가난하다";   English = "1.) adj. poor, indigent, destitute" },
    @{ Foreign = "신혼여행";   English = "1.)  n. honeymoon" },
    @{ Foreign = "경우";       English = "1.) n. case, circumstances, scenario" },
    @{ Foreign = "장난";       English = "1.) n. joke, mischief, prank" },
    @{ Foreign = "장난꾸러기"; English = "1.) mischievous person, jokester, prankster" }
)

$newDate = "2020-12-08"

$startRow = $ws.UsedRange.Rows.Count + 1
$endRow   = $startRow + $newWords.Count - 1

# Format the new DATE cells as Text first, same as a user would via
# Format Cells > Text, so the "yyyy-mm-dd" string is kept verbatim instead
# of being auto-converted into a date serial number.
$dateRangeAddr = "C" + $startRow + ":C" + $endRow
$ws.Range($dateRangeAddr).NumberFormat = "@"

for ($i = 0; $i -lt $newWords.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newWords[$i].Foreign
    $ws.Cells.Item($row, 2).Value = $newWords[$i].English
    $ws.Cells.Item($row, 3).Value = $newDate
}

Write-Host "Added $($newWords.Count) rows ($startRow..$endRow)"
